$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 9 ("gnda" supply-ground row) entirely; rows below shift up by one.
$ws.Rows("9").Delete()

# Former row 11 ("gnd"/"ground") is now row 10; add its comment text.
$ws.Range("H10").Value = "common ground"

# Update the saved selection to match the target state.
$ws.Range("L16").Select()
